$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "TEST"
$ws.Range("J2").Value = "COUPE"

$ws.Range("D4").Select()
